# Apply cryptos list update (GitHub Actions data refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column keeps its original text formatting (avoid Excel
# auto-converting numeric-looking strings like "1.00" or "7.20" into numbers,
# which would silently drop meaningful trailing zeros / dot groupings).
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "57.428.88"
$ws.Range("E2").Value = "  -0.05%  "

$ws.Range("D3").Value = "3.115.08"
$ws.Range("E3").Value = "  +0.48%  "

$ws.Range("E4").Value = "  +0.05%  "

$ws.Range("D5").Value = "525.72"
$ws.Range("E5").Value = "  +0.55%  "

$ws.Range("D6").Value = "137.01"
$ws.Range("E6").Value = "  -2.80%  "

$ws.Range("E7").Value = "  +0.08%  "

$ws.Range("D8").Value = "3.111.51"
$ws.Range("E8").Value = "  +0.38%  "

$ws.Range("D9").Value = "0.446"
$ws.Range("E9").Value = "  +2.54%  "

$ws.Range("D10").Value = "7.20"
$ws.Range("E10").Value = "  -0.50%  "

$ws.Range("E11").Value = "  -0.85%  "

$ws.Range("D12").Value = "0.396"
$ws.Range("E12").Value = "  +3.05%  "

$ws.Range("D13").Value = "3.661.56"
$ws.Range("E13").Value = "  +0.74%  "

$ws.Range("E14").Value = "  +2.63%  "

$ws.Range("D15").Value = "25.32"
$ws.Range("E15").Value = "  -3.03%  "

$ws.Range("E16").Value = "  +0.05%  "

$ws.Range("D17").Value = "57.627.24"
$ws.Range("E17").Value = "  +0.14%  "

$ws.Range("D18").Value = "3.123.72"
$ws.Range("E18").Value = "  +0.65%  "

$ws.Range("D19").Value = "5.95"
$ws.Range("E19").Value = "  -2.51%  "

$ws.Range("D20").Value = "12.60"
$ws.Range("E20").Value = "  -1.38%  "

$ws.Range("D21").Value = "7.89"
$ws.Range("E21").Value = "  -2.01%  "

$ws.Range("D22").Value = "348.07"
$ws.Range("E22").Value = "  +3.59%  "

$ws.Range("D23").Value = "5.79"
$ws.Range("E23").Value = "  -0.53%  "

$ws.Range("D25").Value = "68.26"
$ws.Range("E25").Value = "  +2.66%  "

$ws.Range("D26").Value = "0.504"
$ws.Range("E26").Value = "  -1.37%  "

$ws.Range("D27").Value = "0.167"
$ws.Range("E27").Value = "  -0.61%  "

$ws.Range("D28").Value = "1.01"
$ws.Range("E28").Value = "  +0.36%  "

$ws.Range("D29").Value = "0.0₃0910"
$ws.Range("E29").Value = "  -0.95%  "

$ws.Range("D30").Value = "7.43"
$ws.Range("E30").Value = "  +3.29%  "

$ws.Range("D31").Value = "0.997"
$ws.Range("E31").Value = "  -0.13%  "

$ws.Range("E32").Value = "  +0.38%  "

$ws.Range("D33").Value = "6.06"
$ws.Range("E33").Value = "  -6.75%  "

$ws.Range("D34").Value = "21.05"
$ws.Range("E34").Value = "  +0.82%  "

$ws.Range("D35").Value = "1.17"
$ws.Range("E35").Value = "  -2.05%  "

$ws.Range("D36").Value = "4.96"
$ws.Range("E36").Value = "  +7.01%  "

$ws.Range("D37").Value = "157.96"
$ws.Range("E37").Value = "  +0.67%  "

$ws.Range("D38").Value = "6.15"
$ws.Range("E38").Value = "  +0.64%  "

$ws.Range("D39").Value = "26.27"
$ws.Range("E39").Value = "  -2.67%  "

$ws.Range("E40").Value = "  -2.81%  "

$ws.Range("D41").Value = "0.0667"
$ws.Range("E41").Value = "  +1.04%  "

$ws.Range("D42").Value = "4.17"
$ws.Range("E42").Value = "  +6.07%  "

$ws.Range("D43").Value = "1.61"
$ws.Range("E43").Value = "  +6.95%  "

$ws.Range("D44").Value = "0.699"
$ws.Range("E44").Value = "  +2.11%  "

$ws.Range("D45").Value = "3.160.62"
$ws.Range("E45").Value = "  +0.62%  "

$ws.Range("D46").Value = "36.52"
$ws.Range("E46").Value = "  -0.75%  "

$ws.Range("D47").Value = "0.0269"
$ws.Range("E47").Value = "  +3.81%  "

$ws.Range("B48").Value = "FirstDigitalUSD"
$ws.Range("C48").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D48").Value = "1.00"
$ws.Range("E48").Value = "  +0.07%  "

$ws.Range("B49").Value = "Maker"
$ws.Range("C49").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D49").Value = "2.343.32"
$ws.Range("E49").Value = "  +2.00%  "

$ws.Range("D50").Value = "0.959"
$ws.Range("E50").Value = "  -1.77%  "

$ws.Range("D51").Value = "6.03"
$ws.Range("E51").Value = "  +0.33%  "
